$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7371370792388916
$ws.Range("B1").Value = 2.541946411132812
$ws.Range("C1").Value = 4.940074920654297
$ws.Range("D1").Value = 2.885887145996094
$ws.Range("E1").Value = 0.8921214938163757
